$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-18 Saturday", "2025-10-19 Sunday"),
    @("496÷2=", "741÷6="),
    @("253÷4=", "145÷6="),
    @("283÷6=", "744÷4="),
    @("244÷7=", "281÷4="),
    @("993÷5=", "205÷2="),
    @("577÷3=", "324÷5="),
    @("753÷8=", "973÷5="),
    @("337÷9=", "434÷7="),
    @("440÷8=", "294÷9="),
    @("741÷5=", "394÷3="),
    @("646÷6=", "372÷6="),
    @("131÷2=", "768÷7="),
    @("371÷6=", "477÷5="),
    @("581÷4=", "362÷8="),
    @("340÷8=", "304÷5="),
    @("148÷5=", "866÷4="),
    @("534÷6=", "423÷6="),
    @("225÷7=", "798÷2="),
    @("439÷8=", "138÷6="),
    @("737÷7=", "489÷6="),
    @("993÷3=", "597÷8="),
    @("621÷8=", "119÷3="),
    @("703÷3=", "851÷5="),
    @("566÷4=", "224÷2="),
    @("234÷2=", "887÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
